$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header cell H2: "total" (sharedstring index shuffles automatically) ---
$ws.Range("H2").Value = "total"

# --- Row 3: axes steppers ---
$ws.Range("E3").Value = "axes steppers"
$ws.Range("F3").Value = 4
$ws.Range("G3").Value = 5

# --- Row 4: axes endstops (G4 keeps formula =G3) ---
$ws.Range("E4").Value = "axes endstops"
$ws.Range("F4").Value = 1
$ws.Range("G4").Formula = "=G3"

# Re-create the shared formula group for H3:H4 in one shot (matches the
# author's partial shared-formula range after their edit broke the old
# H3:H10 group).
$ws.Range("H3:H4").Formula = "=F3*G3"

# --- Row 5: axes encoders (now styled like the other data rows, with a
# newly introduced G5 formula referencing G3) ---
$ws.Range("E5").Value = "axes encoders"
$ws.Range("F5").Value = 2
$ws.Range("G5").Formula = "=G3"
$ws.Range("H5").Formula = "=F5*G5"
$ws.Range("E5:F5").Style = "Comma"
$ws.Range("E5").Style = "Normal"
$ws.Range("F5").Style = "Normal"
$ws.Range("E5").Font.Color = 0
$ws.Range("F5").Font.Color = 0

# --- Row 6: door switch ---
$ws.Range("E6").Value = "door switch"
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 1
$ws.Range("H6").Formula = "=F6*G6"

# --- Row 7: tool zeroing (style reset to Normal/default, losing its old
# shaded-font styling) ---
$ws.Range("E7").Value = "tool zeroing"
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 1
$ws.Range("H7").Formula = "=F7*G7"
$ws.Range("E7").Style = "Normal"
$ws.Range("F7").Style = "Normal"
$ws.Range("G7").Style = "Normal"

# --- Row 8: sd card ---
$ws.Range("E8").Value = "sd card"
$ws.Range("F8").Value = 4
$ws.Range("G8").Value = 1
$ws.Range("H8").Formula = "=F8*G8"

# --- Row 9: spindle on/off ---
$ws.Range("E9").Value = "spindle on/off"
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 1
$ws.Range("H9").Formula = "=F9*G9"

# --- Row 10: spindle rpm (style reset to Normal/default) ---
$ws.Range("C10").Value = "debug"
$ws.Range("E10").Value = "spindle rpm"
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 1
$ws.Range("H10").Formula = "=F10*G10"
$ws.Range("E10").Style = "Normal"
$ws.Range("F10").Style = "Normal"
$ws.Range("G10").Style = "Normal"

# --- Row 11: debug port (was the blank "total" summary row, now a normal
# data row) ---
$ws.Range("C11").Value = "debug"
$ws.Range("E11").Value = "debug port"
$ws.Range("F11").Value = 2
$ws.Range("G11").Value = 1
$ws.Range("H11").Formula = "=F11*G11"
$ws.Range("G11").Style = "Normal"
$ws.Range("G11").Font.Color = 0
$ws.Range("H11").Font.Bold = $false
$ws.Range("H11").Font.Color = 0

# --- Row 12: buzzer (was the blank "available" summary row) ---
$ws.Range("E12").Value = "buzzer"
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 1
$ws.Range("H12").Formula = "=F12*G12"
$ws.Range("G12").Style = "Normal"
$ws.Range("G12").Font.Color = 0
$ws.Range("H12").Font.Bold = $false
$ws.Range("H12").Font.Color = 0

# --- Row 13: pwm mosfets (was the blank "remaining" summary row) ---
$ws.Range("E13").Value = "pwm mosfets"
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 3
$ws.Range("H13").Formula = "=F13*G13"
$ws.Range("G13").Style = "Normal"
$ws.Range("G13").Font.Color = 0
$ws.Range("H13").Font.Bold = $false
$ws.Range("H13").Font.Color = 0

# --- Row 14: analog inputs (brand new data row) ---
$ws.Range("E14").Value = "analog inputs"
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 3
$ws.Range("H14").Formula = "=F14*G14"
$ws.Range("E14:H14").Font.Color = 0

# --- Row 15: current monitoring (brand new data row) ---
$ws.Range("E15").Value = "current monitoring"
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 1
$ws.Range("H15").Formula = "=F15*G15"
$ws.Range("E15:H15").Font.Color = 0

# --- Row 16: new "total" summary row (shifted down from row 11) ---
$ws.Range("G16").Value = "total"
$ws.Range("H16").Formula = "=SUM(H3:H15)"
$ws.Range("E16:F16").Font.Color = 0
$ws.Range("G16:H16").Font.Bold = $true
$ws.Range("G16:H16").Font.Color = 0

# --- Row 17: new "available" summary row (shifted down from row 12) ---
$ws.Range("G17").Value = "available"
$ws.Range("H17").Value = 58
$ws.Range("E17:F17").Font.Color = 0
$ws.Range("G17:H17").Font.Bold = $true
$ws.Range("G17:H17").Font.Color = 0

# --- Row 18: new "remaining" summary row (shifted down from row 13) ---
$ws.Range("G18").Value = "remaining"
$ws.Range("H18").Formula = "=H17-H16"
$ws.Range("E18:F18").Font.Color = 0
$ws.Range("G18:H18").Font.Bold = $true
$ws.Range("G18:H18").Font.Color = 0

# --- Move the active selection to G20, as left by the author ---
$ws.Range("G20").Select()
